$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update poll input numbers (row 2)
$ws.Range("A2").Value = 35
$ws.Range("B2").Value = 31
$ws.Range("C2").Value = 11
$ws.Range("D2").Value = 9
$ws.Range("I2").Value = 9

# Update second table input numbers (row 10)
$ws.Range("A10").Value = 49
$ws.Range("B10").Value = 47

# Move selection to K7
$ws.Range("K7").Select()
